# schedule_modules.xlsx update
# Commit: "Adding 10 Manhattan plots lecture and recitation"
#
# Week 11 (Manhattan plots) and Week 12 (Making lots of plots at once) are merged
# into a single "Manhattan plots and making lots of plots at once" session, and the
# remaining weeks of the term are filled in with real topics (ggplot extensions /
# complexheatmap, Thanksgiving break, and the capstone wrap-up), replacing the old
# placeholder "Open for input" rows, on both the "Schedule" and "Schedule_date" sheets.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # "Schedule"
$ws2 = $wb.Worksheets.Item(2)   # "Schedule_date"

# ---------------------------------------------------------------------------
# Final Module / Topic content, by week (1-16), shared by both sheets.
# ---------------------------------------------------------------------------
$modules = @(
    "1: Principles",
    "1: Principles",
    "2: Coding fundamentals",
    "2: Coding fundamentals",
    "2: Coding fundamentals",
    "3: Data exploration",
    "3: Data exploration",
    "Open session, capstone prep",
    "3: Data exploration",
    "4: Putting it together",
    "4: Putting it together",
    "4: Putting it together",
    "4: Putting it together",
    "No class, Thanksgiving",
    "4: Putting it together",
    "4: Putting it together"
)

$topics = @(
    "Principles of data visualization",
    "Good and bad visualizations",
    "R Markdown for reproducible research",
    "ggplot 101",
    "Themes, labels, facets (ggplot 102)",
    "Data distributions",
    "Correlations",
    "Open session, capstone prep",
    "Annotating statistics",
    "Principal components analysis",
    "Manhattan plots and making lots of plots at once",
    "Interactive plots",
    "ggplot extension packages and complexheatmap",
    "Relaxing and eating",
    "Capstone assignment open session",
    "Capstone assignment open session"
)

$dates = @(
    44796, 44803, 44810, 44817, 44825, 44832, 44838, 44845,
    44852, 44859, 44866, 44873, 44880, 44887, 44894, 44900
)

# ---------------------------------------------------------------------------
# Sheet "Schedule": Week | Module | Topic  (rows 2-17, header in row 1)
# ---------------------------------------------------------------------------
for ($i = 0; $i -lt 16; $i++) {
    $r = $i + 2
    $ws1.Cells.Item($r, 1).Value = $i + 1
    $ws1.Cells.Item($r, 2).Value = $modules[$i]
    $ws1.Cells.Item($r, 3).Value = $topics[$i]
}

# ---------------------------------------------------------------------------
# Sheet "Schedule_date": Week | Date | Module | Topic (rows 2-17, header row 1)
# ---------------------------------------------------------------------------
for ($i = 0; $i -lt 16; $i++) {
    $r = $i + 2
    $ws2.Cells.Item($r, 1).Value = $i + 1
    $ws2.Cells.Item($r, 2).Value = $dates[$i]
    $ws2.Cells.Item($r, 3).Value = $modules[$i]
    $ws2.Cells.Item($r, 4).Value = $topics[$i]
}

# ---------------------------------------------------------------------------
# Column widths
# ---------------------------------------------------------------------------
$ws1.Columns.Item(2).ColumnWidth = 24          # raw width -> 24.8333 (was 19.832)
$ws2.Columns.Item(2).ColumnWidth = 12.2        # raw width -> 13
$ws2.Columns.Item(3).ColumnWidth = 21.7        # raw width -> 22.5

# ---------------------------------------------------------------------------
# View / selection state: make "Schedule_date" active first, then "Schedule"
# last so "Schedule" ends up the selected tab, matching the target file.
# ---------------------------------------------------------------------------
$ws2.Activate()
$ws2.Range("A1:D17").Select() | Out-Null

$ws1.Activate()
$ws1.Range("F23").Select() | Out-Null
